$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so Excel does not
# auto-coerce numeric-looking strings (losing trailing zeros, e.g. "1.00" -> 1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.372.88"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.229.45"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -1.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.49"
$ws.Range("E5").Value = "  -2.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.33"
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.560"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.44"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.03"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.564.40"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.230.81"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.38"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.781"
$ws.Range("E17").Value = "  -6.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.135.82"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -5.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.50"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.40"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.24"
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.41"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.34"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.84"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.52"
$ws.Range("E32").Value = "  -7.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0765"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -5.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  -5.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.85"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.69"
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0301"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  -6.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.53"
$ws.Range("E42").Value = "  -9.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.795.33"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.83"
$ws.Range("E45").Value = "  +9.43%  "
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "69.25"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.97"
$ws.Range("E48").Value = "  -4.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.77"
$ws.Range("E49").Value = "  -6.69%  "

# Rows 50/51: FraxShare and THORChain swap order, plus updated price/volume
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("E50").Value = "  -5.79%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -4.49%  "
